# Updated cryptos list on Mon Mar  6 19:36:19 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to be treated as text (so numeric-looking strings like
    # "288.52" are not silently converted to numbers), then restore the
    # cell's original (default/"Normal") style so no visual formatting
    # changes are introduced.
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "22.473.90"
$ws.Range("E2").Value = "  +0.11%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.571.52"
$ws.Range("E3").Value = "  +0.15%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.06%  "

# Row 5 - USDC
Set-TextValue "D5" "1.004"
$ws.Range("E5").Value = "  +0.12%  "

# Row 6 - BNB
Set-TextValue "D6" "288.52"
$ws.Range("E6").Value = "  -0.69%  "

# Row 7 - XRP
Set-TextValue "D7" "0.3735"
$ws.Range("E7").Value = "  +1.17%  "

# Row 8 - OKB
Set-TextValue "D8" "48.30"
$ws.Range("E8").Value = "  -3.00%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.3340"
$ws.Range("E9").Value = "  -0.94%  "

# Row 10 - was Polygon, now Dogecoin (rows 10/11 swapped)
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue "D10" "0.07497"
$ws.Range("E10").Value = "  -0.39%  "

# Row 11 - was Dogecoin, now Polygon
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D11" "1.134"
$ws.Range("E11").Value = "  -1.07%  "

# Row 12 - BinanceUSD
$ws.Range("E12").Value = "  +0.07%  "

# Row 13 - Solana
Set-TextValue "D13" "20.95"
$ws.Range("E13").Value = "  -0.79%  "

# Row 14 - Polkadot
Set-TextValue "D14" "5.982"
$ws.Range("E14").Value = "  -0.52%  "

# Row 15 - Chainlink
Set-TextValue "D15" "6.922"
$ws.Range("E15").Value = "  -0.55%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "1.576.44"
$ws.Range("E16").Value = "  +0.34%  "

# Row 17 - ShibaInu
$ws.Range("E17").Value = "  -0.07%  "

# Row 18 - Litecoin
Set-TextValue "D18" "88.43"
$ws.Range("E18").Value = "  -2.18%  "

# Row 19 - TRON
Set-TextValue "D19" "0.06768"
$ws.Range("E19").Value = "  -0.05%  "

# Row 20 - Dai
Set-TextValue "D20" "1.004"
$ws.Range("E20").Value = "  +0.22%  "

# Row 21 - Uniswap
Set-TextValue "D21" "6.396"
$ws.Range("E21").Value = "  +0.71%  "

# Row 22 - Avalanche
Set-TextValue "D22" "16.48"
$ws.Range("E22").Value = "  +0.64%  "

# Row 23 - Cosmos
Set-TextValue "D23" "12.11"
$ws.Range("E23").Value = "  -0.94%  "

# Row 24 - WrappedBTC (E unchanged)
Set-TextValue "D24" "22.463.13"

# Row 25 - Toncoin
Set-TextValue "D25" "2.391"
$ws.Range("E25").Value = "  +0.42%  "

# Row 26 - LidoDAOToken (E unchanged)
Set-TextValue "D26" "2.580"

# Row 27 - Monero
Set-TextValue "D27" "152.43"
$ws.Range("E27").Value = "  +2.33%  "

# Row 28 - EthereumClassic
Set-TextValue "D28" "19.76"
$ws.Range("E28").Value = "  -1.31%  "

# Row 29 - HuobiToken
Set-TextValue "D29" "5.007"
$ws.Range("E29").Value = "  -0.88%  "

# Row 30 - BitcoinCash
Set-TextValue "D30" "124.27"
$ws.Range("E30").Value = "  -0.59%  "

# Row 31 - WrappedliquidstakedEther2.0
Set-TextValue "D31" "1.749.74"
$ws.Range("E31").Value = "  -0.06%  "

# Row 32 - ImmutableX
Set-TextValue "D32" "1.055"
$ws.Range("E32").Value = "  -0.79%  "

# Row 33 - Filecoin
Set-TextValue "D33" "6.173"
$ws.Range("E33").Value = "  -0.03%  "

# Row 34 - WEMIXTOKEN
Set-TextValue "D34" "2.014"
$ws.Range("E34").Value = "  +0.10%  "

# Row 35 - FraxShare
Set-TextValue "D35" "9.699"
$ws.Range("E35").Value = "  -0.90%  "

# Row 36 - Stellar (E unchanged)
Set-TextValue "D36" "0.08317"

# Row 37 - VeChain
Set-TextValue "D37" "0.02464"
$ws.Range("E37").Value = "  -0.39%  "

# Row 38 - Algorand
Set-TextValue "D38" "0.2276"
$ws.Range("E38").Value = "  -0.95%  "

# Row 39 - Hedera
Set-TextValue "D39" "0.06391"
$ws.Range("E39").Value = "  -2.58%  "

# Row 40 - InternetComputer(DFINITY)
Set-TextValue "D40" "5.393"
$ws.Range("E40").Value = "  -0.64%  "

# Row 41 - TrustWalletToken
Set-TextValue "D41" "1.293"
$ws.Range("E41").Value = "  -4.33%  "

# Row 42 - was TheSandbox, now Aptos (rows 42/43 swapped)
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D42" "11.33"
$ws.Range("E42").Value = "  +1.04%  "

# Row 43 - was Aptos, now TheSandbox
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue "D43" "0.6321"
$ws.Range("E43").Value = "  +1.73%  "

# Row 44 - Frax
$ws.Range("E44").Value = "  +0.09%  "

# Row 45 - EnergySwap
Set-TextValue "D45" "13.93"
$ws.Range("E45").Value = "  -1.51%  "

# Row 46 - Decentraland
Set-TextValue "D46" "0.6150"
$ws.Range("E46").Value = "  +5.18%  "

# Row 47 - PancakeSwap
Set-TextValue "D47" "3.783"
$ws.Range("E47").Value = "  -0.61%  "

# Row 48 - NEARProtocol
Set-TextValue "D48" "2.060"
$ws.Range("E48").Value = "  -0.38%  "

# Row 49 - Quant
Set-TextValue "D49" "125.45"
$ws.Range("E49").Value = "  -1.93%  "

# Row 50 - EOS
Set-TextValue "D50" "1.216"
$ws.Range("E50").Value = "  -1.70%  "

# Row 51 - Cronos
Set-TextValue "D51" "0.07270"
$ws.Range("E51").Value = "  -0.48%  "
